$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (232) down through the new rows (233-240)
$ws.Range("A232:Z232").Copy()
$ws.Range("A233:Z240").PasteSpecial(-4122)

# Row 233
$ws.Range("A233").Value = 45663
$ws.Range("B233").Value = 904.8648037680001
$ws.Range("C233").Value = 255.6743238075
$ws.Range("I233").Value = 390.710493282
$ws.Range("K233").Value = 130.278232933332
$ws.Range("N233").Value = 53.77079296287999
$ws.Range("O233").Value = 1.2418156674
$ws.Range("Q233").Value = 0.0000049872
$ws.Range("U233").Value = 507.2385568937764
$ws.Range("Z233").Value = 1908.311204413872

# Row 234
$ws.Range("A234").Value = 45664
$ws.Range("B234").Value = 858.1239230958001
$ws.Range("C234").Value = 234.4476936185
$ws.Range("I234").Value = 361.973881527
$ws.Range("K234").Value = 122.609097015498
$ws.Range("N234").Value = 45.65686295744
$ws.Range("O234").Value = 1.1883410847
$ws.Range("Q234").Value = 0.0000044616
$ws.Range("U234").Value = 445.4488411359721
$ws.Range("Z234").Value = 1796.983734716812

# Row 235
$ws.Range("A235").Value = 45665
$ws.Range("B235").Value = 841.3605457758001
$ws.Range("C235").Value = 230.7021439915
$ws.Range("I235").Value = 353.576703388
$ws.Range("K235").Value = 121.589781482115
$ws.Range("N235").Value = 42.79312295552
$ws.Range("O235").Value = 1.1849701941
$ws.Range("Q235").Value = 0.000004276799999999999
$ws.Range("U235").Value = 421.3981857271373
$ws.Range("Z235").Value = 1877.288570598708

# Row 236
$ws.Range("A236").Value = 45666
$ws.Range("B236").Value = 819.1617274422001
$ws.Range("C236").Value = 223.20757792
$ws.Range("I236").Value = 331.142451345
$ws.Range("K236").Value = 116.590281485046
$ws.Range("N236").Value = 40.72098588096
$ws.Range("O236").Value = 1.1670091356
$ws.Range("Q236").Value = 0.0000040824
$ws.Range("U236").Value = 399.6502526446804
$ws.Range("Z236").Value = 1787.947111653536

# Row 237
$ws.Range("A237").Value = 45667
$ws.Range("B237").Value = 838.3999598658
$ws.Range("C237").Value = 226.524628904
$ws.Range("I237").Value = 336.030361008
$ws.Range("K237").Value = 118.677451386735
$ws.Range("N237").Value = 41.05858124704
$ws.Range("O237").Value = 1.1801692287
$ws.Range("Q237").Value = 0.0000042336
$ws.Range("U237").Value = 408.4773548957953
$ws.Range("Z237").Value = 1862.429381128682

# Row 238
$ws.Range("A238").Value = 45668
$ws.Range("B238").Value = 837.2836994922001
$ws.Range("C238").Value = 227.6194498705
$ws.Range("I238").Value = 336.764442679
$ws.Range("K238").Value = 116.97859216443
$ws.Range("N238").Value = 40.40667295392
$ws.Range("O238").Value = 1.186110849
$ws.Range("Q238").Value = 0.0000043248
$ws.Range("U238").Value = 404.8953423880964
$ws.Range("Z238").Value = 1813.38008171822

# Row 239
$ws.Range("A239").Value = 45669
$ws.Range("B239").Value = 836.7975261468
$ws.Range("C239").Value = 226.542656355
$ws.Range("I239").Value = 337.516428781
$ws.Range("K239").Value = 113.386718380128
$ws.Range("N239").Value = 39.12613880672
$ws.Range("O239").Value = 1.1803905498
$ws.Range("Q239").Value = 0.00000426
$ws.Range("U239").Value = 396.0682401369816
$ws.Range("Z239").Value = 1756.877897307118

# Row 240
$ws.Range("A240").Value = 45670
$ws.Range("B240").Value = 836.7182231580001
$ws.Range("C240").Value = 217.5434914885
$ws.Range("I240").Value = 327.418329697
$ws.Range("K240").Value = 108.144524208444
$ws.Range("N240").Value = 36.99579563456
$ws.Range("O240").Value = 1.1723889408
$ws.Range("Q240").Value = 0.0000040464
$ws.Range("U240").Value = 385.8339186864135
$ws.Range("Z240").Value = 1786.177051465884

Write-Output "done"
